$wb = $excel.ActiveWorkbook
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsTrans = $wb.Worksheets.Item("Transactions")

# Insert a new (blank) column before column N, shifting the existing
# "Late" / "Outstanding" / "Over Due" columns one place to the right.
$wsRepay.Columns("N").Insert()
$wsRepay.Columns("N").ColumnWidth = $wsRepay.Columns("M").ColumnWidth

# Update the selected range on the "Transactions" sheet and make sure it
# is no longer the active tab.
$wsTrans.Activate()
$wsTrans.Range("D22").Select() | Out-Null

# Update the selected range on the "Repayment schedule" sheet and make
# it the active tab.
$wsRepay.Activate()
$wsRepay.Range("M12").Select() | Out-Null
